$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: lowercase the column headers
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "address"

# Row 4: expand abbreviated street address
$ws.Range("B4").Value = "1940 West 33rd Street, Chicago"

# Row 7: expand abbreviated street address
$ws.Range("B7").Value = "1 South State Street, Chicago"

# Row 13: expand abbreviated street address
$ws.Range("B13").Value = "6525 West Diversey Avenue, Chicago"

# Rows 14/15 swap content (old row 14 becomes old row 15's data, and vice versa)
$ws.Range("A14").Value = "Target"
$ws.Range("B14").Value = "7100 South Cicero Avenue, Bedford Park"

$ws.Range("A15").Value = "Target Market News"
$ws.Range("B15").Value = "228 South Wabash Avenue, Chicago"

# Row 16: replace with expanded address that used to belong to row 17, then delete row 17
$ws.Range("A16").Value = "Target"
$ws.Range("B16").Value = "4050 North Harlem Avenue, Norridge"

# Delete the now-duplicate last row (old row 17)
$ws.Rows.Item(17).Delete()
